$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$orig1 = $ws.Range("L232").NumberFormat
$orig2 = $ws.Range("M232").NumberFormat
$ws.Range("L232").NumberFormat = "General"
$ws.Range("M232").NumberFormat = "General"
$ws.Range("L232").Value = 0
$ws.Range("M232").Value = 0
$ws.Range("L232").NumberFormat = $orig1
$ws.Range("M232").NumberFormat = $orig2
